$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Quantum Mechanics: A Bridge to the Subatomic Realm" "Unraveling the Art of Persuasion in Government and Politics"

# Author name (three runs "Dr" / "." / " Emily Carter" collapse into one)
Replace-Text "Dr. Emily Carter" "Charlotte Robertson"

# Email address (two separate edits, each touching one run but not the "." run)
Replace-Text "ecarter@cambridge" "crobertson@schoolmail"
Replace-Text "edu" "com"

# Body paragraph 1 (sentence by sentence, matching the diff's run layout)
Replace-Text "In the depths of the microscopic universe, where particles dance at the subatomic level, lies a realm governed by the enigmatic laws of quantum mechanics" "The art of persuasion lies at the core of successful government and politics"

Replace-Text " This fascinating branch of physics has revolutionized our understanding of the universe and opened up new avenues of exploration" " In the intricate world of public affairs, the ability to effectively communicate ideas, sway opinions, and influence decisions plays a pivotal role in shaping policy, mobilizing masses, and driving societal change"

Replace-Text " Quantum mechanics is a journey into the realm of the very small, where particles behave in ways that challenge our classical intuition" " As budding citizens, it is essential for us to understand the nuances of persuasion, its strategies and techniques, and the profound impact it has on the functioning of our political systems"

Replace-Text " The study of quantum mechanics has unveiled a world of particles possessing dual characteristics, acting as both waves and particles, as revealed through renowned experiments like the double-slit experiment. This mind-boggling nature of quantum mechanics has led to numerous applications, including the development of cutting-edge technologies like quantum computing" " This exploration into the art of persuasion will unveil the captivating dynamic between leaders, audiences, and the transformation of societies"

Replace-Text "Delving deeper into the principles of quantum mechanics, we encounter concepts such as wave-particle duality, superposition, and entanglement" "Navigating the complex landscape of persuasion requires an understanding of its multifaceted nature"

Replace-Text " These seemingly counterintuitive phenomena have profound implications for our understanding of the physical world" " Persuasion encompasses a spectrum of strategies, from logical reasoning and empathetic appeals to emotional tactics and strategic communication"

Replace-Text " Wave-particle duality blurs the classical distinction between particles and waves, as particles can exhibit properties of both" " In a political context, leaders harness the power of rhetoric, employing skillful oratory and compelling narratives to capture hearts and minds"

Replace-Text " Superposition allows particles to exist in multiple states simultaneously, challenging our conventional notion of reality" " They engage in debates, public speeches, and media appearances, using language as a potent weapon to define issues, sway opinions, and ultimately, mobilize action"

Replace-Text " Entanglement, on the other hand, depicts the inexplicable interconnectedness of particles, even when they are separated by vast distances. These fundamental concepts have ignited debates among physicists and philosophers alike, leading to new insights into the nature of reality itself" " Moreover, the art of persuasion extends beyond mere speech; it manifests in symbolic gestures, policy decisions, and strategic alliances, all working in tandem to paint a persuasive picture of a better future"

Replace-Text "In addition to its theoretical significance, quantum mechanics has also propelled technological advances" "The art of persuasion in government and politics carries immense significance"

Replace-Text " The exploration of quantum effects has paved the way for groundbreaking applications in fields such as computing, cryptography, and communication" " It serves as a catalyst for societal change, shaping public perception, galvanizing movements, and influencing the course of history"

Replace-Text " Quantum computers leverage the quantum properties of particles to perform calculations that are beyond the capacity of classical computers" " In the hands of skilled leaders, persuasion can inspire hope, ignite passion, and bridge seemingly insurmountable divides"

Replace-Text " Quantum cryptography offers unbreakable encryption methods, ensuring the security of sensitive information" " However, it is crucial to recognize the potential for abuse"

Replace-Text " Quantum communication, utilizing entangled particles, facilitates secure and ultra-fast data transmission" " When persuasion becomes manipulation, it can mislead, deceive, and exploit vulnerable populations for personal or political gain"

Replace-Text " These advancements exemplify the practical implications of quantum mechanics and its potential to revolutionize various industries" " Therefore, it is imperative for us to develop a critical understanding of persuasion techniques, discerning genuine attempts to serve the public good from those driven by self-interest or malicious intent"

# Summary heading: move the lastRenderedPageBreak marker here (it used to sit on
# "and ultra-fast data transmission", which was merged away above).
$d.Content.Find.Execute("Summary", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null
if ($d.Content.Find.Found) {
    $d.Content.InsertBefore([char]0)
}

# Summary paragraph text
Replace-Text "Quantum mechanics has reshaped our understanding of the subatomic realm, revealing a world governed by unique principles that challenge our classical intuition" "This essay explores the art of persuasion in the context of government and politics, highlighting its importance in shaping public opinion, influencing policy decisions, and mobilizing societal change"

Replace-Text " Core concepts such as wave-particle duality, superposition, and entanglement have brought forth new insights into the nature of reality" " Through a deep dive into persuasive strategies, ranging from logical reasoning to emotional appeals, the essay illustrates how leaders leverage rhetoric, symbolic gestures, and policy choices to capture hearts and minds"

Replace-Text " The practical applications of quantum mechanics are equally remarkable, leading to innovations in computing, cryptography, communication, and beyond. As we continue to delve into the mysteries of the quantum world, we can anticipate even more transformative applications and a deeper appreciation for the intricate workings of the universe at its most fundamental level" " While recognizing the transformative potential of persuasion, the essay cautions against its misuse and emphasizes the need for critical thinking to discern genuine attempts to serve the public good from those driven by selfish agendas"

# Add a new empty paragraph after the Summary paragraph (before sectPr)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
